$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz2")
$ws.Activate()

# Strip the "Sample_" prefix from the library_name (column B) values for rows 47-58
# (these duplicated the sample directory's "Sample_" prefix already present in
# column D / the localization path, so it is redundant there).
$ws.Range("B47").Value = "GB_RNA_stress_16st_1"
$ws.Range("B48").Value = "GB_RNA_stress_pH5_1"
$ws.Range("B49").Value = "GB_RNA_stress_pH7_1"
$ws.Range("B50").Value = "GB_RNA_stress_pH9_1"
$ws.Range("B51").Value = "GB_RNA_stress_16st_2"
$ws.Range("B52").Value = "GB_RNA_stress_pH5_2"
$ws.Range("B53").Value = "GB_RNA_stress_pH7_2"
$ws.Range("B54").Value = "GB_RNA_stress_pH9_2"
$ws.Range("B55").Value = "GB_RNA_stress_16st_3"
$ws.Range("B56").Value = "GB_RNA_stress_pH5_3"
$ws.Range("B57").Value = "GB_RNA_stress_pH7_3"
$ws.Range("B58").Value = "GB_RNA_stress_pH9_3"

# Update the saved view state: scrolled down so row 40 is the top visible row,
# and the active cell/selection moved to B59 (just past the last data row).
$win = $excel.ActiveWindow
$win.ScrollRow = 40
$win.ScrollColumn = 1
$ws.Range("B59").Select()
